# Complete 100% shoulder and bicept for both arms
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header rename: "Direction" -> "Dir2"
$ws.Range("H1").Value = "Dir2"

# Row 2 - left_shoulder_x / M1
$ws.Range("E2").Value = "Inside"
$ws.Range("F2").Value = 130
$ws.Range("G2").Value = 270
$ws.Range("H2").Value = "Outside"

# Row 3 - left_shoulder_y / M2
$ws.Range("F3").Value = 122

# Row 4 - left_shoulder_z / M3
$ws.Range("D4").Value = 120
$ws.Range("E4").Value = "Inside"
$ws.Range("F4").Value = 123
$ws.Range("G4").Value = 270
$ws.Range("H4").Value = "Outside"

# Row 5 - left_bicept / M4
$ws.Range("D5").Value = 30
$ws.Range("F5").Value = 35
$ws.Range("G5").Value = 125

# Last selected cell
$ws.Range("E2").Select()
